# Swap bonferroni-correction markers with fdr-correction markers:
# Specific "< .05" significance labels become "< .05*" on both the
# "arousal" and "valence" result sheets.

$wb = $excel.ActiveWorkbook

$arousalCells = @("H11","C16","H18","C21","H21","H22","H23","C26","H33")
$valenceCells = @("H11","H13","H14","H15","C16","C20","C21","H23","C26","H28","C31","C32","C40","C41","C42","C43","C44","C45","C48")

$wsArousal = $wb.Worksheets.Item("arousal")
foreach ($addr in $arousalCells) {
    $wsArousal.Range($addr).Value = "< .05*"
}

$wsValence = $wb.Worksheets.Item("valence")
foreach ($addr in $valenceCells) {
    $wsValence.Range($addr).Value = "< .05*"
}
